# Update column F (dSF) values for specific rows, repulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    10 = 1
    15 = 4
    17 = -2
    19 = -1
    20 = 2
    29 = 0
    30 = -2
    31 = 1
    34 = -2
    37 = 3
    38 = -1
    49 = -2
    52 = -1
    53 = 2
    56 = 0
    58 = -3
    63 = 0
    66 = -2
    69 = 0
    74 = 2
    75 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
